$p = $ppt.ActivePresentation

# --- Remove "Talk Outline" (slide 2) and the blank slide (slide 4) ---
# Delete from the highest index down so earlier indices stay stable.
$p.Slides.Item(4).Delete()
$p.Slides.Item(2).Delete()

# --- Slide that now sits at index 2 (was slide3.xml, "405") gets a title
#     and new body copy describing the conda/spack install steps. ---
$install = $p.Slides.Item(2)

$install.Shapes.Item(1).TextFrame.TextRange.Text = "Installation instructions"

$content = $install.Shapes.Item(2).TextFrame.TextRange
$para2 = $content.Paragraphs(2)
$para2.Text = "Before the tutorial, install the "
$para2.InsertAfter("conda") | Out-Null
$para2.InsertAfter(" and ") | Out-Null
$para2.InsertAfter("spack") | Out-Null
$para2.InsertAfter(" components described here:") | Out-Null
$para2.ParagraphFormat.Bullet.Type = 0
